# Update "想去人数" (number of people interested) counts for two events
# that appear on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" -> rows 3 and 4, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1184
$wsExpo.Range("F4").Value = 2650

# Sheet "全部类型" -> rows 5 and 6, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1184
$wsAll.Range("F6").Value = 2650
